$wb = $excel.ActiveWorkbook

# New event row shared by both "展览" (sheet 1) and "全部类型" (sheet 4):
# "南宁·2024良牙动漫秋季盛典（秋典）" starting 2024-10-03.
$newDate  = "2024-10-03"
$newName  = "南宁·2024良牙动漫秋季盛典（秋典）"
$newPlace = "民族大道106号 南宁国际会展中心"
$newRange = "2024.10.03 09:30-10.04 17:30"
$newWant  = 4
$newPrice = 55
$newLink  = "https://show.bilibili.com/platform/detail.html?id=90762"
$newCover = "//i0.hdslb.com/bfs/openplatform/202408/njVhnU591723691579900.jpeg"

function Update-Sheet {
    param($ws)

    # Bump the "want to go" counters on the first three existing events.
    $ws.Cells.Item(2, 6).Value = 1423
    $ws.Cells.Item(3, 6).Value = 2985
    $ws.Cells.Item(4, 6).Value = 32

    # Insert a new row right before the current row 5, shifting every
    # following row (and its data) down by one.
    $ws.Rows.Item(5).Insert()

    # Row 5's numbering cell (col A) should look exactly like the one
    # above it (bold / bordered / centered style).
    $ws.Cells.Item(4, 1).Copy()
    $ws.Cells.Item(5, 1).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Application.CutCopyMode = $false

    $ws.Cells.Item(5, 1).Value = 4

    # Column B holds plain text dates ("2024-08-17", ...) - force text so
    # Excel doesn't silently convert the literal into a date serial, then
    # drop the now-unneeded explicit number format again.
    $ws.Cells.Item(5, 2).NumberFormat = "@"
    $ws.Cells.Item(5, 2).Value = $newDate
    $ws.Cells.Item(5, 2).ClearFormats()

    $ws.Cells.Item(5, 3).Value = $newName
    $ws.Cells.Item(5, 4).Value = $newPlace
    $ws.Cells.Item(5, 5).Value = $newRange
    $ws.Cells.Item(5, 6).Value = $newWant
    $ws.Cells.Item(5, 7).Value = $newPrice
    $ws.Cells.Item(5, 8).Value = $newLink
    $ws.Cells.Item(5, 9).Value = $newCover

    # Column A is simply the 0-based row index (row r -> value r-1); fix
    # up every row pushed down by the insert.
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 6; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

Update-Sheet($wb.Worksheets.Item(1))
Update-Sheet($wb.Worksheets.Item(4))
